$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure the Price column keeps its textual representation (e.g. "5.940", "1.670.64")
# instead of being auto-converted to numbers by Excel, by pre-formatting as Text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.188.48'
$ws.Cells.Item(2, 5).Value = '  -0.07%  '
$ws.Cells.Item(3, 4).Value = '1.670.64'
$ws.Cells.Item(3, 5).Value = '  -0.67%  '
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  -0.27%  '
$ws.Cells.Item(5, 4).Value = '211.05'
$ws.Cells.Item(5, 5).Value = '  -2.41%  '
$ws.Cells.Item(6, 4).Value = '0.5212'
$ws.Cells.Item(6, 5).Value = '  -0.65%  '
$ws.Cells.Item(7, 5).Value = '  -0.25%  '
$ws.Cells.Item(8, 5).Value = '  -2.45%  '
$ws.Cells.Item(9, 4).Value = '0.06327'
$ws.Cells.Item(9, 5).Value = '  -0.49%  '
$ws.Cells.Item(10, 4).Value = '21.22'
$ws.Cells.Item(10, 5).Value = '  -1.05%  '
$ws.Cells.Item(11, 4).Value = '0.07548'
$ws.Cells.Item(11, 5).Value = '  -1.07%  '
$ws.Cells.Item(12, 4).Value = '1.674.10'
$ws.Cells.Item(12, 5).Value = '  -0.21%  '
$ws.Cells.Item(13, 4).Value = '4.445'
$ws.Cells.Item(13, 5).Value = '  -1.57%  '
$ws.Cells.Item(14, 4).Value = '0.5492'
$ws.Cells.Item(14, 5).Value = '  -4.44%  '
$ws.Cells.Item(15, 4).Value = '0.000008034'
$ws.Cells.Item(15, 5).Value = '  -2.69%  '
$ws.Cells.Item(16, 4).Value = '66.45'
$ws.Cells.Item(16, 5).Value = '  +0.54%  '
$ws.Cells.Item(17, 4).Value = '26.194.54'
$ws.Cells.Item(17, 5).Value = '  -0.17%  '
$ws.Cells.Item(18, 5).Value = '  -0.30%  '
$ws.Cells.Item(19, 4).Value = '4.753'
$ws.Cells.Item(19, 5).Value = '  -2.34%  '
$ws.Cells.Item(20, 4).Value = '187.28'
$ws.Cells.Item(20, 5).Value = '  -1.12%  '
$ws.Cells.Item(21, 5).Value = '  -4.12%  '
$ws.Cells.Item(22, 4).Value = '6.214'
$ws.Cells.Item(22, 5).Value = '  -0.41%  '
$ws.Cells.Item(23, 4).Value = '1.004'
$ws.Cells.Item(23, 5).Value = '  -0.23%  '
$ws.Cells.Item(24, 4).Value = '149.94'
$ws.Cells.Item(24, 5).Value = '  +0.98%  '
$ws.Cells.Item(25, 4).Value = '0.1242'
$ws.Cells.Item(25, 5).Value = '  -1.38%  '
$ws.Cells.Item(26, 4).Value = '7.494'
$ws.Cells.Item(26, 5).Value = '  -3.76%  '
$ws.Cells.Item(27, 4).Value = '15.84'
$ws.Cells.Item(27, 5).Value = '  +0.12%  '
$ws.Cells.Item(28, 4).Value = '0.06332'
$ws.Cells.Item(28, 5).Value = '  +0.59%  '
$ws.Cells.Item(29, 4).Value = '1.358'
$ws.Cells.Item(29, 5).Value = '  -1.47%  '
$ws.Cells.Item(30, 4).Value = '1.282'
$ws.Cells.Item(30, 5).Value = '  -2.47%  '
$ws.Cells.Item(31, 4).Value = '3.525'
$ws.Cells.Item(31, 5).Value = '  -1.31%  '
$ws.Cells.Item(32, 4).Value = '3.416'
$ws.Cells.Item(32, 5).Value = '  -4.19%  '
$ws.Cells.Item(33, 4).Value = '1.645'
$ws.Cells.Item(33, 5).Value = '  -2.22%  '
$ws.Cells.Item(34, 4).Value = '1.005'
$ws.Cells.Item(34, 5).Value = '  -1.83%  '
$ws.Cells.Item(35, 4).Value = '0.6053'
$ws.Cells.Item(35, 5).Value = '  -0.96%  '
$ws.Cells.Item(36, 5).Value = '  -0.62%  '
$ws.Cells.Item(37, 4).Value = '2.754'
$ws.Cells.Item(37, 5).Value = '  +0.14%  '
$ws.Cells.Item(38, 4).Value = '1.113.01'
$ws.Cells.Item(38, 5).Value = '  +1.35%  '
$ws.Cells.Item(39, 4).Value = '6.134'
$ws.Cells.Item(39, 5).Value = '  -0.78%  '
$ws.Cells.Item(40, 5).Value = '  +0.14%  '
$ws.Cells.Item(41, 4).Value = '0.8662'
$ws.Cells.Item(41, 5).Value = '  -2.13%  '
$ws.Cells.Item(42, 5).Value = '  -0.50%  '
$ws.Cells.Item(43, 4).Value = '100.42'
$ws.Cells.Item(43, 5).Value = '  -0.03%  '
$ws.Cells.Item(44, 4).Value = '1.824.09'
$ws.Cells.Item(44, 5).Value = '  -0.40%  '
$ws.Cells.Item(45, 4).Value = '0.00000000109'
$ws.Cells.Item(45, 5).Value = '  -0.76%  '
$ws.Cells.Item(46, 4).Value = '55.63'
$ws.Cells.Item(46, 5).Value = '  -3.08%  '
$ws.Cells.Item(47, 5).Value = '  -0.30%  '
$ws.Cells.Item(48, 4).Value = '8.063'
$ws.Cells.Item(48, 5).Value = '  -0.27%  '
$ws.Cells.Item(49, 4).Value = '0.05239'
$ws.Cells.Item(49, 5).Value = '  -0.66%  '
$ws.Cells.Item(50, 4).Value = '0.4243'
$ws.Cells.Item(50, 5).Value = '  -0.83%  '
$ws.Cells.Item(51, 4).Value = '5.940'
$ws.Cells.Item(51, 5).Value = '  -1.19%  '
